$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Add the new "Indexed Matrix" sheet right after Sheet1 ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Indexed Matrix"

# --- Header row (row 2): column index labels 0..5 in B2:G2 ---
$cols = @("B","C","D","E","F","G")
for ($j = 0; $j -lt 6; $j++) {
    $ws2.Range($cols[$j] + "2").Value = $j
}

# --- Matrix body: row index in column A, "row,col" strings in B:G ---
# Shared-string interning order must be: diagonal (0,0)..(5,5) first,
# then each column top-to-bottom skipping the diagonal cell.
for ($i = 0; $i -lt 6; $i++) {
    $ws2.Range("A" + ($i + 3)).Value = $i
}
for ($k = 0; $k -lt 6; $k++) {
    $ws2.Range($cols[$k] + ($k + 3)).Value = "$k,$k"
}
for ($j = 0; $j -lt 6; $j++) {
    for ($i = 0; $i -lt 6; $i++) {
        if ($i -ne $j) {
            $ws2.Range($cols[$j] + ($i + 3)).Value = "$i,$j"
        }
    }
}

# --- Styling ---
# Style "3": thin border all around, centered, default font -> header row A2:G2
$hdrSeed = $ws2.Range("B2")
$hdrSeed.Borders.LineStyle = 1
$hdrSeed.HorizontalAlignment = -4108
$hdrSeed.VerticalAlignment = -4108
$hdrSeed.Copy()
$ws2.Range("A2:G2").PasteSpecial(-4122)

# Style "4": thin border all around, centered, bold red font -> body A3:G8
# Seed the font from Sheet1!B2 (already fontId=1: bold + red) so no new font
# entries are introduced, then layer on the border + centering.
$ws1.Range("B2").Copy()
$bodySeed = $ws2.Range("A3")
$bodySeed.PasteSpecial(-4122)
$bodySeed.Borders.LineStyle = 1
$bodySeed.HorizontalAlignment = -4108
$bodySeed.VerticalAlignment = -4108
$bodySeed.Copy()
$ws2.Range("A3:G8").PasteSpecial(-4122)

# --- Column widths ---
$ws2.Range("A1:G1").ColumnWidth = 8.6

# --- View / selection state ---
$ws2.Range("A2").Select()
$ws1.Range("A11").Select()
$excel.Application.CutCopyMode = $false
